$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns retain their original plain-text formatting
# when values are re-written, regardless of whether they look numeric.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.596.85"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "1.883.98"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "246.31"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "0.4736"
$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "0.2895"
$ws.Range("E8").Value = "  -1.07%  "

$ws.Range("D9").Value = "0.06541"
$ws.Range("E9").Value = "  +0.30%  "

$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("D11").Value = "0.7628"
$ws.Range("E11").Value = "  +2.84%  "

$ws.Range("D12").Value = "99.78"
$ws.Range("E12").Value = "  +2.97%  "

$ws.Range("D13").Value = "0.07823"
$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").Value = "1.882.18"
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("D15").Value = "5.238"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").Value = "284.01"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("D17").Value = "30.570.02"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").Value = "0.000007520"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "2.126.86"
$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").Value = "5.349"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").Value = "6.431"
$ws.Range("E24").Value = "  +2.30%  "

$ws.Range("D25").Value = "9.175"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").Value = "163.74"

$ws.Range("D27").Value = "19.01"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "1.906"
$ws.Range("E28").Value = "  -0.91%  "

$ws.Range("D29").Value = "0.09758"
$ws.Range("E29").Value = "  -0.47%  "

$ws.Range("D30").Value = "1.328"
$ws.Range("E30").Value = "  -1.13%  "

$ws.Range("D31").Value = "1.501"
$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("D32").Value = "4.253"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("D33").Value = "4.189"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").Value = "0.04843"
$ws.Range("E34").Value = "  -1.28%  "

$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  -0.29%  "

$ws.Range("D36").Value = "0.6993"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").Value = "2.773"
$ws.Range("E37").Value = "  +2.37%  "

$ws.Range("D38").Value = "0.01905"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").Value = "2.868"
$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("D40").Value = "6.317"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").Value = "75.57"
$ws.Range("E41").Value = "  -0.76%  "

$ws.Range("D42").Value = "1.978"
$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("D43").Value = "0.4253"
$ws.Range("E43").Value = "  -1.09%  "

$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "0.8384"

$ws.Range("D46").Value = "9.989"
$ws.Range("E46").Value = "  +4.12%  "

$ws.Range("D47").Value = "101.49"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").Value = "7.017"
$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("D50").Value = "0.05777"
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("D51").Value = "0.3958"
$ws.Range("E51").Value = "  -0.58%  "
